$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 13, shifting existing rows 13:125 down to 14:126.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with its data.
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Macroferia Regional de Talca"
$ws.Range("C13").Value = "Maule"
$ws.Range("D13").Value = 45230
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 300000000
$ws.Range("G13").Value = "Espárragos"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 2500
$ws.Range("K13").Value = 1100
$ws.Range("L13").Value = 1100
$ws.Range("M13").Value = 1100
$ws.Range("N13").Value = "`$/kilo"
$ws.Range("O13").Value = "Provincia de Linares"
$ws.Range("P13").Value = 1100
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = "Hortaliza"
